$d = $word.ActiveDocument

# Locate the paragraph containing Hugo Waterfall's name (the author line that
# currently ends with "| " but has no student ID after it yet).
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*Hugo Waterfall*") {
        $target = $para
        break
    }
}

$r = $target.Range
# Exclude the trailing paragraph mark so the new run lands inside this
# paragraph, right after the existing "| " run.
[void]$r.MoveEnd(1, -1)
$r.Collapse(0)

# Insert the missing student ID as a new run.
$r.InsertAfter("40245720")

# Match the Arial formatting used by the rest of the document / paragraph.
$r.Font.Name = "Arial"
$r.Font.NameAscii = "Arial"
$r.Font.NameOther = "Arial"
$r.Font.NameBi = "Arial"

Write-Output "Inserted student ID after Hugo Waterfall's '| '"
